$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.978.11"
$ws.Range("E2").Value = "'  -4.33%  "

$ws.Range("D3").Value = "'2.918.85"
$ws.Range("E3").Value = "'  -7.12%  "

$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'475.85"
$ws.Range("E5").Value = "'  -9.53%  "

$ws.Range("D6").Value = "'128.62"
$ws.Range("E6").Value = "'  -3.29%  "

$ws.Range("E7").Value = "'  -0.10%  "

$ws.Range("D8").Value = "'2.920.16"
$ws.Range("E8").Value = "'  -7.05%  "

$ws.Range("E9").Value = "'  -8.78%  "

$ws.Range("D10").Value = "'6.73"
$ws.Range("E10").Value = "'  -7.13%  "

$ws.Range("D11").Value = "'0.0987"
$ws.Range("E11").Value = "'  -10.79%  "

$ws.Range("D12").Value = "'0.339"
$ws.Range("E12").Value = "'  -12.81%  "

$ws.Range("E13").Value = "'  -2.40%  "

$ws.Range("D14").Value = "'3.418.72"
$ws.Range("E14").Value = "'  -7.15%  "

$ws.Range("D15").Value = "'23.68"
$ws.Range("E15").Value = "'  -8.24%  "

$ws.Range("D16").Value = "'54.864.80"
$ws.Range("E16").Value = "'  -4.58%  "

$ws.Range("D17").Value = "'2.913.83"
$ws.Range("E17").Value = "'  -7.28%  "

$ws.Range("D18").Value = "'0.0000136"
$ws.Range("E18").Value = "'  -11.10%  "

$ws.Range("D19").Value = "'5.48"
$ws.Range("E19").Value = "'  -5.75%  "

$ws.Range("D20").Value = "'11.64"
$ws.Range("E20").Value = "'  -11.38%  "

$ws.Range("D21").Value = "'7.21"
$ws.Range("E21").Value = "'  -10.37%  "

$ws.Range("D22").Value = "'305.48"
$ws.Range("E22").Value = "'  -11.89%  "

$ws.Range("E23").Value = "'  +0.04%  "

$ws.Range("D24").Value = "'0.449"
$ws.Range("E24").Value = "'  -11.98%  "

$ws.Range("D25").Value = "'59.50"
$ws.Range("E25").Value = "'  -14.27%  "

$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "'  -0.40%  "

$ws.Range("D27").Value = "'0.154"
$ws.Range("E27").Value = "'  -7.18%  "

$ws.Range("E28").Value = "'  +0.04%  "

$ws.Range("D29").Value = "'0.0₃0818"
$ws.Range("E29").Value = "'  -14.39%  "

$ws.Range("D30").Value = "'6.25"
$ws.Range("E30").Value = "'  -8.78%  "

$ws.Range("E31").Value = "'  -8.48%  "

$ws.Range("E32").Value = "'  -6.19%  "

$ws.Range("E33").Value = "'  -12.33%  "

$ws.Range("D34").Value = "'18.99"
$ws.Range("E34").Value = "'  -12.46%  "

$ws.Range("D35").Value = "'145.23"
$ws.Range("E35").Value = "'  -8.87%  "

$ws.Range("E36").Value = "'  -13.80%  "

$ws.Range("D37").Value = "'5.48"
$ws.Range("E37").Value = "'  -12.34%  "

$ws.Range("E38").Value = "'  -12.66%  "

$ws.Range("D39").Value = "'23.43"
$ws.Range("E39").Value = "'  -10.07%  "

$ws.Range("D40").Value = "'0.0627"
$ws.Range("E40").Value = "'  -9.82%  "

$ws.Range("D41").Value = "'2.945.86"
$ws.Range("E41").Value = "'  -7.10%  "

$ws.Range("E42").Value = "'  +0.03%  "

$ws.Range("D43").Value = "'35.69"
$ws.Range("E43").Value = "'  -11.70%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.617"
$ws.Range("E44").Value = "'  -10.69%  "

$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "'0.968"
$ws.Range("E45").Value = "'  -10.38%  "

$ws.Range("E46").Value = "'  -8.34%  "

$ws.Range("E47").Value = "'  -12.58%  "

$ws.Range("D48").Value = "'2.066.40"
$ws.Range("E48").Value = "'  -8.65%  "

$ws.Range("D49").Value = "'5.43"
$ws.Range("E49").Value = "'  -12.77%  "

$ws.Range("E50").Value = "'  -6.22%  "

$ws.Range("E51").Value = "'  -11.53%  "
